# Comentando parte que salva pre-processadas em pasta
# -> Adiciona as estatisticas (knn/tree/svm) das imagens Originais e
#    Preprocessadas em duas tabelas na Planilha1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.7369791666667
$ws.Columns.Item(2).ColumnWidth = 9.16666666666667

# ---------------------------------------------------------------------
# Table 1 - "Original" (A2:D10)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "imagem"
$ws.Range("B2").Value = "knn"
$ws.Range("C2").Value = "tree"
$ws.Range("D2").Value = "svm"

$ws.Range("A3").Value = "04_R_N_H_P"
$ws.Range("B3").Value = 63.1578
$ws.Range("C3").Value = 63.1578
$ws.Range("D3").Value = 63.1578

$ws.Range("A4").Value = "08_R_N_L_P"
$ws.Range("B4").Value = 63.1578
$ws.Range("C4").Value = 63.1578
$ws.Range("D4").Value = 63.1578

$ws.Range("B3:B10").NumberFormat = "0.0000"
$ws.Range("C3:D4").NumberFormat = "0.0000"

$lo1 = $ws.ListObjects.Add(1, $ws.Range("A2:D10"), 0, 1)
$lo1.Name = "Tabela3"
$lo1.TableStyle = "TableStyleMedium1"

# ---------------------------------------------------------------------
# Table 2 - "Preprocessada" (A13:D20)
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "imagem"
$ws.Range("B13").Value = "knn"
$ws.Range("C13").Value = "tree"
$ws.Range("D13").Value = "svm"

$ws.Range("A14").Value = "08_R_N_L_P"
$ws.Range("B14").Value = 60.5263
$ws.Range("C14").Value = 63.1578
$ws.Range("D14").Value = 63.1578

$ws.Range("B14:B20").NumberFormat = "0.0000"
$ws.Range("C14:D14").NumberFormat = "0.0000"

$lo2 = $ws.ListObjects.Add(1, $ws.Range("A13:D20"), 0, 1)
$lo2.Name = "Tabela35"
$lo2.TableStyle = "TableStyleMedium1"

# ---------------------------------------------------------------------
# Title banners - merge first, then style, so every merged cell shares
# the same cell style (matches how Excel records a merged+styled band).
# ---------------------------------------------------------------------
$title1 = $ws.Range("A1:D1")
$title1.Merge()
$ws.Range("A1").Value = "Original"
$title1.Style = "Check Cell"
$title1.HorizontalAlignment = -4108
$title1.Borders.LineStyle = -4119

$title2 = $ws.Range("A12:D12")
$title2.Merge()
$ws.Range("A12").Value = "Preprocessada"
$title2.Style = "Check Cell"
$title2.HorizontalAlignment = -4108
$title2.Borders.LineStyle = -4119

# ---------------------------------------------------------------------
# Row heights to mirror the thick/double-border banner rows
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 15.75

# ---------------------------------------------------------------------
# Selection left wherever the author's cursor was on save
# ---------------------------------------------------------------------
$ws.Range("G16").Select()

Write-Output "done"
